$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hours dedicated on rows 4 and 5 (D column, time-formatted values)
$ws.Range("D4").Value = 4/24
$ws.Range("D5").Value = 3/24

# Recalculate formulas so F4 (=SUM(D:D)) reflects updated totals
$excel.Calculate()

# Update the active cell selection to E12
$ws.Range("E12").Select()
